# Update "horarios" workbook with the latest scrape timestamp and
# recalculated minute counts for Línea 141.

$wb = $excel.ActiveWorkbook

$oldTime = "02:23:01"
$newTime = "02:41:48"

# --- Sheet 1: LP1912 ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $newTime"

$ws1.Range("A6").Value = $newTime
$ws1.Range("D6").Value = 17

$ws1.Range("A7").Value = $newTime
$ws1.Range("D7").Value = 67

$ws1.Range("A8").Value = $newTime
$ws1.Range("D8").Value = 80

# --- Sheet 2: LP1912-215 ------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"

$ws2.Range("A6").Value = $newTime
$ws2.Range("D6").Value = 17

# --- Sheet 3: 6203-6173 --------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
